$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "/home" redirect destination in B4 becomes a real hyperlink to
# https://www.google.com - both the displayed text and the link target.
$ws.Range("B4").Value = "https://www.google.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.google.com") | Out-Null

# Column B is widened to fit the new (longer) URL text.
$ws.Columns.Item(2).ColumnWidth = 22.33203125

# The active selection ends up on B14 after the edit.
$ws.Range("B14").Select() | Out-Null
